$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row labels (drop the "T1" suffix)
$ws.Range("A1").Value = "square"
$ws.Range("B1").Value = "loc1"
$ws.Range("C1").Value = "loc2"
$ws.Range("D1").Value = "corrAns"

# Move the active selection from D1 to D2
$ws.Range("D2").Select()
